$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "At war start she is tied up in San Diego," -> "... on patrol in
#    the Mediterranean."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "At war start she is tied up in San Diego,", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "At war start she is on patrol in the Mediterranean.", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) TG Iowa list: add an air-wing sub-bullet under BB-61 USS Iowa, move
#    CG-23 USS Halsey down below it, and add a helo sub-bullet under
#    FFG-9 USS Wadsworth.
#
#    Before:                        After:
#      BB-61 USS Iowa                 BB-61 USS Iowa
#      CG-23 USS Halsey               VC-6 'Firebees' 5x RQ-2A   (ilvl 1)
#      FFG-9 USS Wadsworth            CG-23 USS Halsey
#                                      FFG-9 USS Wadsworth
#                                      HSL-38 2x SH-2F            (ilvl 1)
# ---------------------------------------------------------------------

# Locate the original "CG-23 USS Halsey" paragraph.
$halsey = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "CG-23 USS Halsey*") {
        $halsey = $cand
        break
    }
}

# First clone it (while still clean/unmodified) into a new paragraph right
# after itself -- this becomes the "new" CG-23 USS Halsey bullet.
$halsey.Range.InsertParagraphAfter() | Out-Null
$cgHalsey = $d.Paragraphs.Item($halsey.Index + 1)
$cgHalsey.Range.ListFormat.ListLevelNumber = 1
$cgHalsey.Range.Text = "CG-23 USS Halsey"

# Now repurpose the original paragraph into the VC-6 sub-bullet (one level
# deeper).
$halsey.Range.ListFormat.ListLevelNumber = 2
$halsey.Range.Font.Size = 11
$halsey.Range.Font.SizeBi = 11
$halsey.Range.Text = "VC-6 'Firebees' 5x RQ-2A"
$halsey.Range.NoProofing = 0

# Locate "FFG-9 USS Wadsworth" paragraph and append the HSL-38 sub-bullet
# right after it.
$wadsworth = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "FFG-9 USS Wadsworth*") {
        $wadsworth = $cand
        break
    }
}

$wadsworth.Range.InsertParagraphAfter() | Out-Null
$hsl = $d.Paragraphs.Item($wadsworth.Index + 1)
$hsl.Range.ListFormat.ListLevelNumber = 2
$hsl.Range.Text = "HSL-38 2x SH-2F"
